$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("run_1")
$ws.Range("F2").Value = 29.96174049377441
$ws.Range("F3").Value = 29.78438615798951
$ws.Range("F4").Value = 29.62546849250793
$ws.Range("F5").Value = 29.62471342086792
$ws.Range("F6").Value = 29.54766011238098
$ws.Range("F7").Value = 29.63524317741394
$ws.Range("F8").Value = 29.60845232009888
$ws.Range("F9").Value = 29.63120079040528
$ws.Range("F10").Value = 29.63078212738037
$ws.Range("F11").Value = 29.94678854942322
$ws.Range("F12").Value = 29.61541700363159
$ws.Range("F13").Value = 29.76275873184204
$ws.Range("F14").Value = 29.57445549964905
$ws.Range("F15").Value = 29.67438411712646
$ws.Range("F16").Value = 29.73203134536743
$ws.Range("F17").Value = 29.76456689834595
$ws.Range("F18").Value = 29.77015519142151
$ws.Range("F19").Value = 29.69072222709656
$ws.Range("F20").Value = 29.57468318939209
$ws.Range("F21").Value = 29.88923645019531

$ws = $wb.Worksheets.Item("run_2")
$ws.Range("F2").Value = 29.98364663124084
$ws.Range("F3").Value = 29.72603487968445
$ws.Range("F4").Value = 29.66450953483581
$ws.Range("F5").Value = 29.80825686454773
$ws.Range("F6").Value = 29.63968539237976
$ws.Range("F7").Value = 29.6733717918396
$ws.Range("F8").Value = 29.68173623085022
$ws.Range("F9").Value = 29.6186830997467
$ws.Range("F10").Value = 29.64041996002197
$ws.Range("F11").Value = 29.92515659332276
$ws.Range("F12").Value = 29.71010065078736
$ws.Range("F13").Value = 29.80621671676636
$ws.Range("F14").Value = 29.55991768836975
$ws.Range("F15").Value = 29.69466066360474
$ws.Range("F16").Value = 29.65447378158569
$ws.Range("F17").Value = 29.46091103553772
$ws.Range("F18").Value = 29.5955741405487
$ws.Range("F19").Value = 29.71961712837219
$ws.Range("F20").Value = 29.80578541755676
$ws.Range("F21").Value = 30.13820433616639

$ws = $wb.Worksheets.Item("run_3")
$ws.Range("F2").Value = 30.08850598335266
$ws.Range("F3").Value = 29.91855692863464
$ws.Range("F4").Value = 29.72693347930908
$ws.Range("F5").Value = 29.87942266464233
$ws.Range("F6").Value = 29.70138549804688
$ws.Range("F7").Value = 29.82395458221436
$ws.Range("F8").Value = 29.86751890182495
$ws.Range("F9").Value = 29.75716161727905
$ws.Range("F10").Value = 29.7739098072052
$ws.Range("F11").Value = 29.95717191696167
$ws.Range("F12").Value = 29.71072292327881
$ws.Range("F13").Value = 29.6784679889679
$ws.Range("F14").Value = 29.78381419181824
$ws.Range("F15").Value = 29.98187208175659
$ws.Range("F16").Value = 29.86907625198364
$ws.Range("F17").Value = 29.81101512908936
$ws.Range("F18").Value = 29.89785242080688
$ws.Range("F19").Value = 29.8581383228302
$ws.Range("F20").Value = 29.88987064361572
$ws.Range("F21").Value = 30.04284644126892

$ws = $wb.Worksheets.Item("run_4")
$ws.Range("F2").Value = 29.98294520378113
$ws.Range("F3").Value = 29.89377474784851
$ws.Range("F4").Value = 29.67446756362915
$ws.Range("F5").Value = 29.75350451469421
$ws.Range("F6").Value = 29.75709795951844
$ws.Range("F7").Value = 29.79875063896179
$ws.Range("F8").Value = 29.73088812828064
$ws.Range("F9").Value = 29.74415707588196
$ws.Range("F10").Value = 29.84790682792664
$ws.Range("F11").Value = 29.98173403739929
$ws.Range("F12").Value = 29.76532912254333
$ws.Range("F13").Value = 29.67596244812012
$ws.Range("F14").Value = 29.6386616230011
$ws.Range("F15").Value = 29.66742587089539
$ws.Range("F16").Value = 29.7730803489685
$ws.Range("F17").Value = 29.71079230308533
$ws.Range("F18").Value = 29.8144896030426
$ws.Range("F19").Value = 29.72393012046814
$ws.Range("F20").Value = 29.85459518432617
$ws.Range("F21").Value = 30.02745604515076

$ws = $wb.Worksheets.Item("run_5")
$ws.Range("F2").Value = 30.02897477149964
$ws.Range("F3").Value = 29.67903447151184
$ws.Range("F4").Value = 29.69495511054993
$ws.Range("F5").Value = 29.93721175193786
$ws.Range("F6").Value = 29.70258188247681
$ws.Range("F7").Value = 29.64367437362671
$ws.Range("F8").Value = 29.80858516693115
$ws.Range("F9").Value = 29.64525747299194
$ws.Range("F10").Value = 29.81223130226136
$ws.Range("F11").Value = 29.90772199630737
$ws.Range("F12").Value = 29.73106503486633
$ws.Range("F13").Value = 29.74583292007446
$ws.Range("F14").Value = 29.7270565032959
$ws.Range("F15").Value = 29.70574951171875
$ws.Range("F16").Value = 29.705162525177
$ws.Range("F17").Value = 29.6982216835022
$ws.Range("F18").Value = 29.7186450958252
$ws.Range("F19").Value = 29.64789414405823
$ws.Range("F20").Value = 29.73085403442383
$ws.Range("F21").Value = 30.00629210472107

Write-Host "Done updating Epoch Time values."